# edit.ps1
# Applies the "added particles to the teleportation" change:
#   1. "Make the normal spell..." paragraph: drop the trailing
#      " and release it when you release the trigger" clause, keeping
#      the final period, and leave the _GoBack bookmark collapsed right
#      before that period.
#   2. 11/26/2016 paragraph: append a new sentence about the particle
#      trail / spell throwing behaviour.
#   3. Remove the _GoBack bookmark from its old location in the
#      11/11/2016 paragraph (merging the two runs it used to separate).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# Change 2: append the new sentence to the end of the 11/26/2016
# paragraph. Doing this first (as a clean "insert at paragraph end")
# creates a brand new trailing run without disturbing the two runs
# that are already there.
# ---------------------------------------------------------------------
$oldParticleText = "11/26/2016 You have to create a particle effect for the platform that is the target of the teleportation. Make the particle effect disappear if the pointer leaves the platform. Make the seal active if you press the grip."
$newSentence = " Shows a trail of particles when teleporting. Throws spells that wait until you release the trigger to move."

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("11/26/2016")) {
        $p.Range.InsertAfter($newSentence)
    }
}

# The two original runs of that paragraph ("11/26/2016 " and "You have
# to create ...") should collapse into a single run, same as the rest
# of the document's runs get re-flowed by Word after an edit, while the
# freshly-appended sentence must remain its own run. A temporary
# bookmark dropped exactly on the boundary between old/new text acts as
# a wall that keeps the upcoming find/replace from merging across it.
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("11/26/2016")) {
        $wallPoint = $p.Range.Start + $oldParticleText.Length
        $wallRange = $d.Range($wallPoint, $wallPoint)
        $d.Bookmarks.Add("TEMP_WALL", $wallRange)
    }
}

$d.Content.Find.Execute($oldParticleText, $true, $false, $false, $false, $false,
                         $true, 1, $false, $oldParticleText, 2)

$d.Bookmarks("TEMP_WALL").Delete()

# ---------------------------------------------------------------------
# Change 1: trim "Make the normal spell to grow with the pressed
# trigger and release it when you release the trigger." down to
# "Make the normal spell to grow with the pressed trigger."
# ---------------------------------------------------------------------
$clause = " and release it when you release the trigger"

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("Make the normal spell")) {
        $text = $p.Range.Text
        $idx = $text.IndexOf($clause)
        $start = $p.Range.Start + $idx
        $delRange = $d.Range($start, $start + $clause.Length)
        $delRange.Delete()
    }
}

# ---------------------------------------------------------------------
# Change 3: move the _GoBack bookmark from the 11/11/2016 paragraph
# (where it separated " that will send it to" from " the player
# class...") to the new edit location, right before the final period of
# the spell paragraph. Removing it merges the two runs it used to sit
# between.
# ---------------------------------------------------------------------
$d.Bookmarks("_GoBack").Delete()

$d.Content.Find.Execute(" that will send it to the player class, which will pass it to the character. ",
                         $true, $false, $false, $false, $false,
                         $true, 1, $false,
                         " that will send it to the player class, which will pass it to the character. ",
                         2)

for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.StartsWith("Make the normal spell")) {
        $text = $p.Range.Text
        $dotIdx = $text.IndexOf(".")
        $bmPoint = $p.Range.Start + $dotIdx
        $bmRange = $d.Range($bmPoint, $bmPoint)
        $d.Bookmarks.Add("_GoBack", $bmRange)
    }
}
